$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.798.27'
$ws.Range("E2").Value = '  -0.90%  '
$ws.Range("D3").Value = '2.033.22'
$ws.Range("E3").Value = '  -1.42%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.91'
$ws.Range("E5").Value = '  -1.49%  '
$ws.Range("E6").Value = '  -0.72%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '59.91'
$ws.Range("E7").Value = '  +3.03%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  -0.21%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0817'
$ws.Range("E10").Value = '  +0.89%  '
$ws.Range("E11").Value = '  +0.61%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.60'
$ws.Range("E12").Value = '  -0.09%  '
$ws.Range("D13").Value = '2.334.40'
$ws.Range("E13").Value = '  -1.46%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.99'
$ws.Range("E14").Value = '  +1.15%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.758'
$ws.Range("E15").Value = '  +0.41%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.21'
$ws.Range("E16").Value = '  -1.48%  '
$ws.Range("D17").Value = '2.035.71'
$ws.Range("E17").Value = '  -1.33%  '
$ws.Range("D18").Value = '37.761.26'
$ws.Range("E18").Value = '  -0.80%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.06'
$ws.Range("E19").Value = '  -2.27%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.78'
$ws.Range("E20").Value = '  -0.10%  '
$ws.Range("D21").Value = '0.0₃0822'
$ws.Range("E21").Value = '  -1.24%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '225.49'
$ws.Range("E22").Value = '  +0.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.39'
$ws.Range("E24").Value = '  -2.46%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.20'
$ws.Range("E25").Value = '  -2.30%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.24'
$ws.Range("E26").Value = '  -0.37%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '164.87'
$ws.Range("E27").Value = '  -0.47%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.128'
$ws.Range("E28").Value = '  -4.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.89'
$ws.Range("E29").Value = '  -0.96%  '
$ws.Range("E30").Value = '  -7.15%  '
$ws.Range("E31").Value = '  +1.12%  '
$ws.Range("E33").Value = '  +3.72%  '
$ws.Range("E34").Value = '  -2.06%  '
$ws.Range("E35").Value = '  -2.97%  '
$ws.Range("E36").Value = '  +6.31%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.24'
$ws.Range("E37").Value = '  -5.73%  '
$ws.Range("E38").Value = '  -1.95%  '
$ws.Range("E39").Value = '  +0.01%  '
$ws.Range("D40").Value = '1.539.65'
$ws.Range("E40").Value = '  +4.05%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '16.94'
$ws.Range("E41").Value = '  +0.66%  '
$ws.Range("E42").Value = '  -1.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '96.75'
$ws.Range("E43").Value = '  -1.75%  '
$ws.Range("E44").Value = '  -1.60%  '
$ws.Range("E45").Value = '  -2.75%  '
$ws.Range("E46").Value = '  -1.68%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.92'
$ws.Range("E47").Value = '  -4.14%  '
$ws.Range("E48").Value = '  -2.06%  '
$ws.Range("E49").Value = '  -0.38%  '
$ws.Range("E50").Value = '  +0.31%  '
$ws.Range("D51").Value = '2.223.76'
$ws.Range("E51").Value = '  -1.44%  '
